# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that each hold a single big table) get
#    their table style switched from the custom "Table_0" style
#    ({C646A5DA-D9C0-4A39-B04E-B0100D949C58}) to the built-in
#    "Medium Style 2 - Accent 1" style ({1BC0F446-9211-484D-8F6F-9644B2E9F55A}).
#
# 2) The slide master's theme color scheme is switched from the deck's
#    custom "Integral / Red Violet" palette to the stock Office palette
#    (what ends up persisted as the "Office Theme" colours in
#    ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{1BC0F446-9211-484D-8F6F-9644B2E9F55A}"
$tableSlideNumbers = @(14, 15, 16)

foreach ($slideNum in $tableSlideNumbers) {
    $slide = $p.Slides.Item($slideNum)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Slide master theme colours -----------------------------------
function Set-ThemeColor($scheme, $index, $hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # RGBColor.RGB is packed like a Win32 COLORREF (0x00BBGGRR)
    $packed = ($bb * 65536) + ($gg * 256) + $rr
    $scheme.Colors($index).RGB = $packed
}

$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$masterColorScheme = $p.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorIndex = $i + 1
    $colorHex = $officeThemeColors[$i]
    Set-ThemeColor $masterColorScheme $colorIndex $colorHex
}
